$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATE_TYPE_CODE (J2): 001 -> 004
# Force the cell to stay text (preserve the leading zero like the original
# inline string "001"), then restore the default "Normal" style so no
# stray number-format style is left attached to the cell.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "004"
$ws.Range("J2").Style = "Normal"

# Update REPORT_DATE (N2): 2018-12-31 -> 2020-09-30
$ws.Range("N2").Value = "2020-09-30 00:00:00"

# Update numeric financial figures for row 2 (O2:AG2)
$ws.Range("O2").Value = 10708070724.45
$ws.Range("P2").Value = 836167016.42
$ws.Range("Q2").Value = 411613225.37
$ws.Range("R2").Value = -70.789347124
$ws.Range("S2").Value = 949697532.5700001
$ws.Range("T2").Value = 17.2104094041
$ws.Range("U2").Value = 395258446.6
$ws.Range("V2").Value = -69.91376571329999
$ws.Range("W2").Value = 6141098938.44
$ws.Range("X2").Value = 1236802943.98
$ws.Range("Y2").Value = 89.80755761109999
$ws.Range("Z2").Value = 3109372.26
$ws.Range("AA2").Value = -99.82194648310001
$ws.Range("AB2").Value = 4566971786.01
$ws.Range("AC2").Value = -48.6060735496
$ws.Range("AD2").Value = -35.5852954015
$ws.Range("AE2").Value = -20.6313254915
$ws.Range("AF2").Value = 70.63586912
$ws.Range("AG2").Value = 57.3501903048
